$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F header (reuse the bold/bordered header style from A1:E1) ---
$ws.Cells.Item(1, 6).Value = "Trening"
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)

# --- Replace the text timestamps in column A (rows 2-7) with real date values ---
$serials = @(
    45684.59146018518,
    45684.59153773148,
    45684.59161180555,
    45684.59145439815,
    45684.59157476852,
    45684.59160833334
)
for ($i = 0; $i -lt $serials.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $serials[$i]
    $ws.Cells.Item($row, 6).Value = "Duża Gra"
}

# --- New rows 8 and 9 (date-only serials, some blank cells, training split) ---
$ws.Cells.Item(8, 1).Value = 45684
$ws.Cells.Item(8, 5).Value = "10-15"
$ws.Cells.Item(8, 6).Value = "Mała Gra"

$ws.Cells.Item(9, 1).Value = 45684
$ws.Cells.Item(9, 5).Value = "5-10"
$ws.Cells.Item(9, 6).Value = "Mała Gra"

# --- Apply date/time number format to all of column A's data cells ---
# First apply the lowercase variant, then switch to the uppercase variant that
# is actually kept on the cells (mirrors the author changing their mind on
# capitalisation, leaving an orphan numFmt behind).
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
